$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FILESTREAM")

# Unify the three lab-results data-dictionary names into a single shared value
$ws.Range("B7").Value = "NHANES-LAB-RESULTS"
$ws.Range("B8").Value = "NHANES-LAB-RESULTS"
$ws.Range("B9").Value = "NHANES-LAB-RESULTS"

# Update the active cell selection on the sheet
$ws.Activate()
$ws.Range("B7").Select()
